$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Round coordinate values in Q and R columns for rows 10 and 11
$ws.Range("Q10").Value = 489385
$ws.Range("R10").Value = 7032365
$ws.Range("Q11").Value = 489462
$ws.Range("R11").Value = 7032627

# Clear the Starttid (Z) and Sluttid (AB) cells for rows 10 and 11
$ws.Range("Z10").ClearContents()
$ws.Range("AB10").ClearContents()
$ws.Range("Z11").ClearContents()
$ws.Range("AB11").ClearContents()
